$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H49").Value = 132.85715
$ws.Range("I49").Value = 133.33333
$ws.Range("J49").Value = 130
$ws.Range("K49").Value = 399.99999
$ws.Range("L49").Value = 390
$ws.Range("M49").Value = -263.99999
$ws.Range("N49").Value = -662
$ws.Range("H58").Value = 482.26315
$ws.Range("I58").Value = 303.70587
$ws.Range("K58").Value = 911.11761
$ws.Range("M58").Value = -761.11761
$ws.Range("H64").Value = 10077.444
$ws.Range("I64").Value = 7199.5
$ws.Range("K64").Value = 7199.5
$ws.Range("M64").Value = -6951.5
$ws.Range("H67").Value = 10077.444
$ws.Range("I67").Value = 7199.5
$ws.Range("K67").Value = 7199.5
$ws.Range("M67").Value = -6341.5
$ws.Range("H92").Value = 6547.647
$ws.Range("I92").Value = 9085.166999999999
$ws.Range("K92").Value = 9085.166999999999
$ws.Range("M92").Value = -7837.166999999999
$ws.Range("H99").Value = 2209.3076
$ws.Range("J99").Value = 4356.1665
$ws.Range("L99").Value = 13068.4995
$ws.Range("N99").Value = -16064.4995
$ws.Range("H106").Value = 1810
$ws.Range("I106").Value = 1113.3334
$ws.Range("K106").Value = 1113.3334
$ws.Range("M106").Value = -482.3334
$ws.Range("H112").Value = 126746.625
$ws.Range("J112").Value = 144640.28
$ws.Range("L112").Value = 433920.84
$ws.Range("N112").Value = -436136.84
$ws.Range("H132").Value = 30170.715
$ws.Range("I132").Value = 34836.035
$ws.Range("K132").Value = 104508.105
$ws.Range("M132").Value = -101978.105
$ws.Range("H137").Value = 2689.4
$ws.Range("I137").Value = 3586
$ws.Range("J137").Value = 2091.6667
$ws.Range("K137").Value = 10758
$ws.Range("L137").Value = 6275.000100000001
$ws.Range("M137").Value = -8208
$ws.Range("N137").Value = -11375.0001
$ws.Range("H138").Value = 3680.7532
$ws.Range("I138").Value = 1598.6666
$ws.Range("J138").Value = 4184.484
$ws.Range("K138").Value = 4795.9998
$ws.Range("L138").Value = 12553.452
$ws.Range("M138").Value = 344.0002000000004
$ws.Range("N138").Value = -22833.452

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 353.94116
$ws.Range("I5").Value = 213.09091
$ws.Range("K5").Value = 213.09091
$ws.Range("M5").Value = -101.09091
$ws.Range("H61").Value = 3513.8948
$ws.Range("I61").Value = 3000.1667
$ws.Range("K61").Value = 3000.1667
$ws.Range("M61").Value = -2788.1667
$ws.Range("H122").Value = 2105.353
$ws.Range("I122").Value = 1723.5834
$ws.Range("K122").Value = 5170.7502
$ws.Range("M122").Value = -2720.7502
$ws.Range("H127").Value = 0
$ws.Range("J127").Value = 0
$ws.Range("L127").Value = 0
$ws.Range("N127").Value = ""
$ws.Range("H129").Value = 50780
$ws.Range("J129").Value = 50780
$ws.Range("L129").Value = 50780
$ws.Range("N129").Value = -60780
$ws.Range("H132").Value = 1930.8334
$ws.Range("I132").Value = 1542.2858
$ws.Range("J132").Value = 3290.75
$ws.Range("K132").Value = 4626.857400000001
$ws.Range("L132").Value = 9872.25
$ws.Range("M132").Value = -2096.857400000001
$ws.Range("N132").Value = -14932.25
$ws.Range("H136").Value = 3513.8948
$ws.Range("I136").Value = 3000.1667
$ws.Range("K136").Value = 9000.500100000001
$ws.Range("M136").Value = -6450.500100000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 353.94116
$ws.Range("I4").Value = 213.09091
$ws.Range("K4").Value = 213.09091
$ws.Range("M4").Value = -98.09091000000001
$ws.Range("H94").Value = 1796.8948
$ws.Range("I94").Value = 2193.2727
$ws.Range("J94").Value = 1251.875
$ws.Range("K94").Value = 2193.2727
$ws.Range("L94").Value = 1251.875
$ws.Range("M94").Value = -1742.2727
$ws.Range("N94").Value = -2153.875

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 1307.4117
$ws.Range("I7").Value = 448.2857
$ws.Range("K7").Value = 448.2857
$ws.Range("M7").Value = -335.2857
$ws.Range("H98").Value = 100777.5
$ws.Range("J98").Value = 100777.5
$ws.Range("L98").Value = 100777.5
$ws.Range("N98").Value = -105269.5
$ws.Range("H112").Value = 66798.2
$ws.Range("J112").Value = 66798.2
$ws.Range("L112").Value = 66798.2
$ws.Range("N112").Value = -69752.2

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 12229.143
$ws.Range("I70").Value = 9001.5
$ws.Range("J70").Value = 13520.2
$ws.Range("K70").Value = 9001.5
$ws.Range("L70").Value = 13520.2
$ws.Range("M70").Value = -8731.5
$ws.Range("N70").Value = -14060.2
$ws.Range("H73").Value = 12229.143
$ws.Range("I73").Value = 9001.5
$ws.Range("J73").Value = 13520.2
$ws.Range("K73").Value = 9001.5
$ws.Range("L73").Value = 13520.2
$ws.Range("M73").Value = -8065.5
$ws.Range("N73").Value = -15392.2
$ws.Range("H97").Value = 1200.4546
$ws.Range("I97").Value = 1017.5
$ws.Range("K97").Value = 1017.5
$ws.Range("M97").Value = -521.5
$ws.Range("H122").Value = 3145.6667
$ws.Range("I122").Value = 2274.5
$ws.Range("J122").Value = 4888
$ws.Range("K122").Value = 6823.5
$ws.Range("L122").Value = 14664
$ws.Range("M122").Value = -4373.5
$ws.Range("N122").Value = -19564
$ws.Range("H132").Value = 235259.77
$ws.Range("J132").Value = 3064.9092
$ws.Range("L132").Value = 9194.7276
$ws.Range("N132").Value = -14254.7276

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 5034.9033
$ws.Range("J16").Value = 6154.5654
$ws.Range("L16").Value = 6154.5654
$ws.Range("N16").Value = -6494.5654
$ws.Range("H42").Value = 12727.75
$ws.Range("I42").Value = 8995
$ws.Range("J42").Value = 18326.875
$ws.Range("K42").Value = 8995
$ws.Range("L42").Value = 18326.875
$ws.Range("M42").Value = -8432
$ws.Range("N42").Value = -19452.875
$ws.Range("H49").Value = 12727.75
$ws.Range("I49").Value = 8995
$ws.Range("J49").Value = 18326.875
$ws.Range("K49").Value = 8995
$ws.Range("L49").Value = 18326.875
$ws.Range("M49").Value = -8848
$ws.Range("N49").Value = -18620.875
$ws.Range("H135").Value = 113332.664
$ws.Range("J135").Value = 113332.664
$ws.Range("L135").Value = 113332.664
$ws.Range("N135").Value = -123472.664
$ws.Range("H136").Value = 6183.476
$ws.Range("I136").Value = 5989.25
$ws.Range("J136").Value = 6442.4443
$ws.Range("K136").Value = 17967.75
$ws.Range("L136").Value = 19327.3329
$ws.Range("M136").Value = -15417.75
$ws.Range("N136").Value = -24427.3329

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H112").Value = 29998.4
$ws.Range("J112").Value = 29998.4
$ws.Range("L112").Value = 29998.4
$ws.Range("N112").Value = -32952.4
$ws.Range("H122").Value = 2995.2727
$ws.Range("I122").Value = 2927.182
$ws.Range("J122").Value = 3131.4546
$ws.Range("K122").Value = 8781.545999999998
$ws.Range("L122").Value = 9394.363799999999
$ws.Range("M122").Value = -6331.545999999998
$ws.Range("N122").Value = -14294.3638
$ws.Range("H127").Value = 0
$ws.Range("J127").Value = 0
$ws.Range("L127").Value = 0
$ws.Range("N127").Value = ""
$ws.Range("H132").Value = 4080.0454
$ws.Range("I132").Value = 3819.1177
$ws.Range("K132").Value = 11457.3531
$ws.Range("M132").Value = -8927.3531
$ws.Range("H136").Value = 1683.72
$ws.Range("I136").Value = 1486.7368
$ws.Range("K136").Value = 4460.2104
$ws.Range("M136").Value = -1910.2104
